$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 28139.379
$ws.Range("J129").Value = 43078.418
$ws.Range("L129").Value = 129235.254
$ws.Range("N129").Value = -139235.254
$ws.Range("H132").Value = 2552517.2
$ws.Range("I132").Value = 2696878.5
$ws.Range("J132").Value = 2133.3333
$ws.Range("K132").Value = 8090635.5
$ws.Range("L132").Value = 6399.999899999999
$ws.Range("M132").Value = -8088105.5
$ws.Range("N132").Value = -11459.9999
$ws.Range("H137").Value = 1236.9231
$ws.Range("I137").Value = 1147.909
$ws.Range("J137").Value = 1726.5
$ws.Range("K137").Value = 3443.727
$ws.Range("L137").Value = 5179.5
$ws.Range("M137").Value = -893.7270000000003
$ws.Range("N137").Value = -10279.5
$ws.Range("H138").Value = 5129682
$ws.Range("I138").Value = 7937215.5
$ws.Range("J138").Value = 2882.6086
$ws.Range("K138").Value = 23811646.5
$ws.Range("L138").Value = 8647.825800000001
$ws.Range("M138").Value = -23806506.5
$ws.Range("N138").Value = -18927.8258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 66
$ws.Range("H74").Value = 608.0909
$ws.Range("I74").Value = 628.2143
$ws.Range("J74").Value = 495.4
$ws.Range("K74").Value = 628.2143
$ws.Range("L74").Value = 495.4
$ws.Range("M74").Value = 245.7857
$ws.Range("N74").Value = -2243.4
$ws.Range("H77").Value = 608.0909
$ws.Range("I77").Value = 628.2143
$ws.Range("J77").Value = 495.4
$ws.Range("K77").Value = 3141.0715
$ws.Range("L77").Value = 2477
$ws.Range("M77").Value = 1226.9285
$ws.Range("N77").Value = -11213
$ws.Range("H102").Value = 1777.8
$ws.Range("I102").Value = 1632.6666
$ws.Range("J102").Value = 1995.5
$ws.Range("K102").Value = 1632.6666
$ws.Range("L102").Value = 1995.5
$ws.Range("M102").Value = -10.66660000000002
$ws.Range("N102").Value = -5239.5
$ws.Range("H132").Value = 23961.8
$ws.Range("I132").Value = 100012
$ws.Range("J132").Value = 4949.25
$ws.Range("K132").Value = 300036
$ws.Range("L132").Value = 14847.75
$ws.Range("M132").Value = -297506
$ws.Range("N132").Value = -19907.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 337.07693
$ws.Range("I64").Value = 268
$ws.Range("K64").Value = 268
$ws.Range("M64").Value = -43
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 337.07693
$ws.Range("I67").Value = 268
$ws.Range("K67").Value = 268
$ws.Range("M67").Value = 512
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 10000
$ws.Range("K76").Value = 10000
$ws.Range("M76").Value = -9685
$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 10000
$ws.Range("K79").Value = 10000
$ws.Range("M79").Value = -8908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1198.5714
$ws.Range("I16").Value = 978
$ws.Range("J16").Value = 1750
$ws.Range("K16").Value = 978
$ws.Range("L16").Value = 1750
$ws.Range("M16").Value = -691
$ws.Range("N16").Value = -2324
$ws.Range("H31").Value = 5559010.5
$ws.Range("I31").Value = 2715
$ws.Range("J31").Value = 10420770
$ws.Range("K31").Value = 2715
$ws.Range("L31").Value = 10420770
$ws.Range("M31").Value = -2420
$ws.Range("N31").Value = -10421360
$ws.Range("H34").Value = 5559010.5
$ws.Range("I34").Value = 2715
$ws.Range("J34").Value = 10420770
$ws.Range("K34").Value = 2715
$ws.Range("L34").Value = 10420770
$ws.Range("M34").Value = -2513
$ws.Range("N34").Value = -10421174
$ws.Range("H113").Value = 1198.5714
$ws.Range("I113").Value = 978
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 978
$ws.Range("L113").Value = 1750
$ws.Range("M113").Value = 1192
$ws.Range("N113").Value = -6090
$ws.Range("H132").Value = 1796.0851
$ws.Range("I132").Value = 1428.5814
$ws.Range("J132").Value = 5746.75
$ws.Range("K132").Value = 4285.7442
$ws.Range("L132").Value = 17240.25
$ws.Range("M132").Value = -1755.7442
$ws.Range("N132").Value = -22300.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 2666.6667
$ws.Range("J93").Value = 2666.6667
$ws.Range("L93").Value = 8000.000100000001
$ws.Range("N93").Value = -11744.0001
$ws.Range("H101").Value = 25500
$ws.Range("J101").Value = 25500
$ws.Range("L101").Value = 76500
$ws.Range("N101").Value = -81368
$ws.Range("H116").Value = 35500
$ws.Range("I116").Value = 50750
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 152250
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = -148808
$ws.Range("N116").Value = -21884
$ws.Range("H131").Value = 342413.5
$ws.Range("I131").Value = 4338.893
$ws.Range("J131").Value = 567796.6
$ws.Range("K131").Value = 13016.679
$ws.Range("L131").Value = 1703389.8
$ws.Range("M131").Value = -7976.679
$ws.Range("N131").Value = -1713469.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2866.75
$ws.Range("I80").Value = 3555
$ws.Range("J80").Value = 2637.3333
$ws.Range("K80").Value = 3555
$ws.Range("L80").Value = 2637.3333
$ws.Range("M80").Value = -2557
$ws.Range("N80").Value = -4633.3333
$ws.Range("H83").Value = 2866.75
$ws.Range("I83").Value = 3555
$ws.Range("J83").Value = 2637.3333
$ws.Range("K83").Value = 17775
$ws.Range("L83").Value = 13186.6665
$ws.Range("M83").Value = -12783
$ws.Range("N83").Value = -23170.6665
$ws.Range("H122").Value = 2196.3726
$ws.Range("I122").Value = 2591.1177
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 7773.353099999999
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -5323.353099999999
$ws.Range("N122").Value = -10897
$ws.Range("H126").Value = 1770.08
$ws.Range("I126").Value = 2219.182
$ws.Range("J126").Value = 1417.2142
$ws.Range("K126").Value = 6657.545999999999
$ws.Range("L126").Value = 4251.642599999999
$ws.Range("M126").Value = -4187.545999999999
$ws.Range("N126").Value = -9191.642599999999
$ws.Range("H132").Value = 22148.203
$ws.Range("I132").Value = 26455.574
$ws.Range("J132").Value = 3004.3333
$ws.Range("K132").Value = 79366.72200000001
$ws.Range("L132").Value = 9012.999899999999
$ws.Range("M132").Value = -76836.72200000001
$ws.Range("N132").Value = -14072.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1866.4667
$ws.Range("I81").Value = 1858.25
$ws.Range("J81").Value = 1899.3334
$ws.Range("K81").Value = 3716.5
$ws.Range("L81").Value = 3798.6668
$ws.Range("M81").Value = -2655.5
$ws.Range("N81").Value = -5920.6668
$ws.Range("H84").Value = 1866.4667
$ws.Range("I84").Value = 1858.25
$ws.Range("J84").Value = 1899.3334
$ws.Range("K84").Value = 18582.5
$ws.Range("L84").Value = 18993.334
$ws.Range("M84").Value = -13278.5
$ws.Range("N84").Value = -29601.334
$ws.Range("H100").Value = 298.625
$ws.Range("I100").Value = 227
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 454
$ws.Range("L100").Value = 1600
$ws.Range("M100").Value = 87
$ws.Range("N100").Value = -2682
$ws.Range("H122").Value = 30078.03
$ws.Range("I122").Value = 33717.805
$ws.Range("K122").Value = 101153.415
$ws.Range("M122").Value = -98703.41500000001
$ws.Range("H132").Value = 1797.25
$ws.Range("I132").Value = 1540.238
$ws.Range("J132").Value = 2568.2856
$ws.Range("K132").Value = 4620.714
$ws.Range("L132").Value = 7704.8568
$ws.Range("M132").Value = -2090.714
$ws.Range("N132").Value = -12764.8568
